{"js": "// Locate the paragraph that ends the DELETE Command bullet list:\n// \"The user must specify the row they wish to delete via the Primary Key\"\n// and insert three new bulleted paragraphs right after it (before the\n// \"ORDER BY Command\" paragraph), describing the new INSERT Command section.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"The user must specify the row they wish to delete via the Primary Key\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph for the INSERT Command insertion point\");\n}\n\nconst newItems = [\n  { level: 0, text: \"INSERT Command\" },\n  { level: 1, text: \"The user must enter values for all the columns that exist in the table\" },\n  { level: 1, text: \"The user must specify the values for the columns in syntax similar to the CREATE TABLE command.\" }\n];\n\nlet previous = anchor;\nfor (const item of newItems) {\n  const p = previous.insertParagraph(item.text, \"After\");\n  // New paragraphs inherit the numbering/list attachment of the paragraph\n  // they were inserted after; drop that inherited attachment before\n  // re-attaching at the level this bullet actually needs.\n  p.detachFromList();\n  await context.sync();\n\n  p.font.set({ name: \"Consolas\", size: 10 });\n  p.attachToList(2, item.level);\n  await context.sync();\n\n  previous = p;\n}\n", "ps1": "# Locate the paragraph that ends the DELETE Command bullet list:\n# \"The user must specify the row they wish to delete via the Primary Key\"\n# and insert three new bulleted paragraphs right after it (before the\n# \"ORDER BY Command\" paragraph), describing the new INSERT Command section.\n$d = $word.ActiveDocument\n\n$anchorText = \"The user must specify the row they wish to delete via the Primary Key\"\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq $anchorText) {\n    $anchorIndex = $i\n    break\n  }\n}\nif ($anchorIndex -eq -1) {\n  throw \"Could not find anchor paragraph for the INSERT Command insertion point\"\n}\n\n$anchor = $d.Paragraphs.Item($anchorIndex)\n$anchor.Range.InsertParagraphAfter()\n\n# New paragraphs inherit the anchor's numbering/level (ilvl=1, numId=2);\n# only the ilvl needs correcting per bullet (ListLevelNumber is 1-based:\n# 1 => ilvl 0, 2 => ilvl 1). numId/style/font all carry over already.\n$p1 = $d.Paragraphs.Item($anchorIndex + 1)\n$p1.Range.Text = \"INSERT Command\"\n$p1.Range.ListFormat.ListLevelNumber = 1\n\n$p1.Range.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Item($anchorIndex + 2)\n$p2.Range.Text = \"The user must enter values for all the columns that exist in the table\"\n$p2.Range.ListFormat.ListLevelNumber = 2\n\n$p2.Range.InsertParagraphAfter()\n$p3 = $d.Paragraphs.Item($anchorIndex + 3)\n$p3.Range.Text = \"The user must specify the values for the columns in syntax similar to the CREATE TABLE command.\"\n$p3.Range.ListFormat.ListLevelNumber = 2\n"}
